$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the materials cost for "edX Course Certificates" (row 27) from the
# Optional column (D) to the Required column (C).
$ws.Range("C27").Value2 = $ws.Range("D27").Value2
$ws.Range("D27").Clear()

# Update the active selection to match the new state (C36).
$ws.Range("C36").Select()
